$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting the existing data (rows 108-138)
# down to rows 109-139, mirroring the weekly refresh that prepends a new
# price observation for Jengibre at Vega Modelo de Temuco.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new weekly record. The
# descriptive columns repeat the same boilerplate values used throughout
# this market/category block.
$ws.Cells.Item(108, 1).Value = 10
$ws.Cells.Item(108, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value = "La Araucanía"
$ws.Cells.Item(108, 4).Value = 44627
$ws.Cells.Item(108, 5).Value = 9
$ws.Cells.Item(108, 6).Value = 100114007
$ws.Cells.Item(108, 7).Value = "Jengibre"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 35
$ws.Cells.Item(108, 11).Value = 25000
$ws.Cells.Item(108, 12).Value = 26000
$ws.Cells.Item(108, 13).Value = 25571
$ws.Cells.Item(108, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(108, 15).Value = "Perú"
$ws.Cells.Item(108, 16).Value = 1967
$ws.Cells.Item(108, 17).Value = 13
$ws.Cells.Item(108, 18).Value = "Hortaliza"
